$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "63.722.23"
Set-TextValue "E2" "  -3.12%  "

Set-TextValue "D3" "2.608.09"
Set-TextValue "E3" "  -2.18%  "

Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.05%  "

Set-TextValue "D5" "573.40"
Set-TextValue "E5" "  -3.98%  "

Set-TextValue "D6" "155.45"
Set-TextValue "E6" "  -1.17%  "

Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.05%  "

Set-TextValue "E8" "  -6.16%  "

Set-TextValue "D9" "0.118"
Set-TextValue "E9" "  -6.65%  "

Set-TextValue "D10" "5.81"
Set-TextValue "E10" "  -0.38%  "

Set-TextValue "D11" "0.380"
Set-TextValue "E11" "  -4.94%  "

Set-TextValue "E12" "  -0.38%  "

Set-TextValue "D13" "28.06"
Set-TextValue "E13" "  -2.55%  "

Set-TextValue "D14" "3.081.11"
Set-TextValue "E14" "  -1.96%  "

Set-TextValue "D15" "0.0000179"
Set-TextValue "E15" "  -7.71%  "

Set-TextValue "D16" "63.516.30"
Set-TextValue "E16" "  -3.22%  "

Set-TextValue "D17" "2.579.11"
Set-TextValue "E17" "  -2.15%  "

Set-TextValue "D18" "12.01"
Set-TextValue "E18" "  -4.47%  "

Set-TextValue "E19" "  +2.52%  "

Set-TextValue "E20" "  -5.06%  "

Set-TextValue "E21" "  -2.14%  "

Set-TextValue "E22" "  -0.09%  "

Set-TextValue "D23" "67.16"
Set-TextValue "E23" "  -3.59%  "

Set-TextValue "E24" "  -1.79%  "

Set-TextValue "E25" "  -3.32%  "

Set-TextValue "D26" "587.21"
Set-TextValue "E26" "  +4.12%  "

Set-TextValue "D27" "9.14"
Set-TextValue "E27" "  -4.30%  "

Set-TextValue "E28" "  -2.68%  "

Set-TextValue "E29" "  +0.00%  "

Set-TextValue "E30" "  -1.51%  "

Set-TextValue "D31" "7.89"

Set-TextValue "E32" "  -2.38%  "

Set-TextValue "E33" "  -3.48%  "

Set-TextValue "D34" "6.51"
Set-TextValue "E34" "  -1.11%  "

Set-TextValue "E35" "  -2.19%  "

Set-TextValue "D36" "0.405"
Set-TextValue "E36" "  -3.95%  "

Set-TextValue "D37" "1.00"
Set-TextValue "E37" "  +0.06%  "

Set-TextValue "D38" "19.62"
Set-TextValue "E38" "  -4.31%  "

Set-TextValue "D39" "154.08"
Set-TextValue "E39" "  -0.27%  "

Set-TextValue "E40" "  -3.95%  "

Set-TextValue "E41" "  -0.04%  "

Set-TextValue "D42" "41.48"
Set-TextValue "E42" "  -2.86%  "

Set-TextValue "E43" "  +5.42%  "

Set-TextValue "D44" "155.55"
Set-TextValue "E44" "  -3.26%  "

Set-TextValue "D45" "3.89"
Set-TextValue "E45" "  -4.58%  "

Set-TextValue "D46" "23.17"
Set-TextValue "E46" "  +2.10%  "

Set-TextValue "D47" "0.0588"
Set-TextValue "E47" "  -2.57%  "

Set-TextValue "D48" "0.627"
Set-TextValue "E48" "  -1.88%  "

Set-TextValue "E49" "  -2.35%  "

Set-TextValue "D50" "0.0245"
Set-TextValue "E50" "  -3.81%  "

Set-TextValue "D51" "18.85"
Set-TextValue "E51" "  -4.75%  "
